# "Updated my Sprint Backlog hours"
# Update the "Amount Remaining After…" hours for a few tasks on the
# "Sprint 1" sheet (rows 11, 13, 14) and let the totals / chart recalc.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sprint 1")

# Row 11 - Myles Debro / Implment UI binding
$ws.Range("D11").Value = 2
$ws.Range("E11").Value = 0.5
$ws.Range("F11").Value = 2
$ws.Range("G11").Value = 1

# Row 13 - Myles Debro / Implment Login feature (admin login idea)
$ws.Range("D13").Value = 0
$ws.Range("E13").Value = 2
$ws.Range("F13").Value = 1
# G13 unchanged (stays 3)

# Row 14 - Myles Debro / Implment Logout feature
$ws.Range("D14").Value = 0
$ws.Range("E14").Value = 0
$ws.Range("F14").Value = 2
$ws.Range("G14").Value = 0.5

# Reflect the author's final cursor position / zoom level in the sheet view
[void]$ws.Range("E11").Select()
$excel.ActiveWindow.Zoom = 117
